# Adapt column header formatting to respective input file names.
# The "_old" / "_new" header suffixes become "_FV2210" / "_FV2304", a
# native Excel Table ("Table1") is laid over the data range so the new
# headers act as the table's column names, and the top header row is
# frozen.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Base header labels (without suffix), in left-to-right column order.
$labels = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J (1-10): "<label>_old" -> "<label>_FV2210"
for ($i = 0; $i -lt $labels.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = "$($labels[$i])_FV2210"
}

# Column K (11) stays "diff" - untouched.

# Columns L-U (12-21): "<label>_new" -> "<label>_FV2304"
for ($i = 0; $i -lt $labels.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = "$($labels[$i])_FV2304"
}

# Turn the data range into a native Excel Table using the (now renamed)
# header row as the column names.
$dataRange = $ws.Range("A1:U57")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# Freeze the header row (View > Freeze Panes > Freeze Top Row).
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Renamed headers, added Table1 over A1:U57, froze header row."
